# Recompute the RANDBETWEEN-driven test-vector block (rows 3-32) and
# stamp column F ("o1") with the running sample index, copied from
# column A ("Time"), for each of the 30 test lines.
#
# - Columns B:D hold volatile =RANDBETWEEN(-5,5) shared formulas; simply
#   touching/recalculating the sheet re-draws new random cached values.
# - Column F holds plain literal values (no formula) that should mirror
#   the Time value (column A) for each row instead of staying 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 32; $row++) {
    $time = $ws.Cells.Item($row, 1).Value2
    $ws.Cells.Item($row, 6).Value = $time
}

# Force a full recalculation so the volatile RANDBETWEEN formulas in
# B:D redraw fresh cached values alongside the column F edits above.
$ws.Calculate()
